$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new column before column N (14th column), shifting the
# "Late" / "Outstanding" columns one place to the right.
$ws.Columns("N").Insert()

# Match the column width Excel used for the freshly inserted column
# (same character width as column M, but without "best fit").
$ws.Columns("N").ColumnWidth = 9.17

# Update the active selection on the sheet to match the saved state.
$ws.Range("R7").Select()
